$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B numeric / value updates ---
$ws.Range("B1").Value = 0.03299177083333334
$ws.Range("B2").Value = 28.10595194444444
$ws.Range("B3").Value = 1416.031751931667
$ws.Range("B5").Value = 10.454
$ws.Range("B6").Value = 98
$ws.Range("B7").Value = 26

$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("B8").Value = 31.43142949901265

$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("B9").Value = 45.05145882646376

$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("B10").Value = 72

$ws.Range("B11").Value = "Custom mode`n98.87%`nEco mode`n0.30%`nSports mode`n0.09%"

$ws.Range("A12").Value = "Peak Power(kW)"

$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("B13").Value = -1800.676194614624

$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B14").Value = 0.04504690499999999

$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.003181106069741899

$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.325

$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3

$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("B18").Value = 0.3250000000000002

$ws.Range("A19").Value = "Minimum Temperature(C)"

$ws.Range("A20").Value = "Maximum Temperature(C)"

$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 14

$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"

$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"

$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"

$ws.Range("A25").Value = "Maximum MCU Temperature(C)"

$ws.Range("A26").Value = "Maximum Motor Temperature(C)"

$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

$ws.Range("A28").Value = "highest cell temp(C)"

$ws.Range("A29").Value = "lowest cell temp(C)"

$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 53

$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.489615453055555

$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.0000001451866913309508

$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 20.40043089038297

$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 2.841647784257643

$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 4.825229374837487

$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 10.63482040043089

$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 10.37480034174065

$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 8.012332379926452

$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 10.14078228891943

$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 13.66962594257271

$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 18.9257457003826

# --- New row 43 ---
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
